$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price column (D) stays text, matching the source data formatting
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '27.697.57'
$ws.Range("E2").Value = '  +0.89%  '
$ws.Range("D3").Value = '1.778.65'
$ws.Range("E3").Value = '  +1.80%  '
$ws.Range("D4").Value = '0.9997'
$ws.Range("E4").Value = '  -0.43%  '
$ws.Range("D5").Value = '326.18'
$ws.Range("E5").Value = '  +0.63%  '
$ws.Range("D6").Value = '0.9986'
$ws.Range("E6").Value = '  -0.39%  '
$ws.Range("D7").Value = '0.4599'
$ws.Range("E7").Value = '  +4.06%  '
$ws.Range("D8").Value = '0.3589'
$ws.Range("E8").Value = '  -0.49%  '
$ws.Range("D9").Value = '0.07487'
$ws.Range("E9").Value = '  +0.44%  '
$ws.Range("D10").Value = '41.94'
$ws.Range("E10").Value = '  -0.40%  '
$ws.Range("D11").Value = '1.105'
$ws.Range("E11").Value = '  +0.85%  '
$ws.Range("D12").Value = '0.9986'
$ws.Range("E12").Value = '  -0.42%  '
$ws.Range("D13").Value = '20.90'
$ws.Range("E13").Value = '  +1.70%  '
$ws.Range("D14").Value = '6.050'
$ws.Range("E14").Value = '  +0.65%  '
$ws.Range("D15").Value = '7.222'
$ws.Range("E15").Value = '  +1.42%  '
$ws.Range("D16").Value = '1.770.05'
$ws.Range("E16").Value = '  +1.09%  '
$ws.Range("D17").Value = '93.78'
$ws.Range("E17").Value = '  +2.21%  '
$ws.Range("D18").Value = '0.00001058'
$ws.Range("E18").Value = '  +0.20%  '
$ws.Range("D19").Value = '0.06420'
$ws.Range("E19").Value = '  +0.42%  '
$ws.Range("D20").Value = '0.9985'
$ws.Range("E20").Value = '  -0.29%  '
$ws.Range("D21").Value = '17.11'
$ws.Range("E21").Value = '  +1.86%  '
$ws.Range("D22").Value = '5.797'
$ws.Range("E22").Value = '  -1.00%  '
$ws.Range("D23").Value = '27.762.91'
$ws.Range("E23").Value = '  +0.92%  '
$ws.Range("D24").Value = '11.31'
$ws.Range("E24").Value = '  +1.16%  '
$ws.Range("D25").Value = '2.082'
$ws.Range("E25").Value = '  -1.10%  '
$ws.Range("D26").Value = '164.96'
$ws.Range("E26").Value = '  +2.25%  '
$ws.Range("D27").Value = '20.35'
$ws.Range("E27").Value = '  +0.09%  '
$ws.Range("D28").Value = '1.978.37'
$ws.Range("E28").Value = '  +1.36%  '
$ws.Range("D29").Value = '2.169'
$ws.Range("E29").Value = '  +4.40%  '
$ws.Range("D30").Value = '126.39'
$ws.Range("E30").Value = '  +1.67%  '
$ws.Range("D31").Value = '1.102'
$ws.Range("E31").Value = '  +2.72%  '
$ws.Range("D32").Value = '0.09221'
$ws.Range("E32").Value = '  +3.01%  '
$ws.Range("D33").Value = '3.666'
$ws.Range("E33").Value = '  +0.37%  '
$ws.Range("D34").Value = '5.565'
$ws.Range("E34").Value = '  +1.46%  '
$ws.Range("D35").Value = '11.83'
$ws.Range("E35").Value = '  -0.93%  '
$ws.Range("D36").Value = '0.02296'
$ws.Range("E36").Value = '  -0.02%  '
$ws.Range("D37").Value = '0.06110'
$ws.Range("E37").Value = '  +2.45%  '
$ws.Range("D38").Value = '0.2093'
$ws.Range("E38").Value = '  +0.66%  '
$ws.Range("D39").Value = '0.6316'
$ws.Range("E39").Value = '  -0.02%  '
$ws.Range("D40").Value = '4.968'
$ws.Range("E40").Value = '  +0.89%  '
$ws.Range("D41").Value = '1.183'
$ws.Range("E41").Value = '  -1.25%  '
$ws.Range("D42").Value = '1.394'
$ws.Range("E42").Value = '  +0.74%  '
$ws.Range("D43").Value = '7.816'
$ws.Range("E43").Value = '  +1.04%  '
$ws.Range("D44").Value = '13.25'
$ws.Range("E44").Value = '  +0.20%  '
$ws.Range("D45").Value = '3.729'
$ws.Range("E45").Value = '  +0.56%  '
$ws.Range("D46").Value = '0.5899'
$ws.Range("E46").Value = '  +0.54%  '
$ws.Range("D47").Value = '122.63'
$ws.Range("E47").Value = '  +1.26%  '
$ws.Range("D48").Value = '1.953'
$ws.Range("E48").Value = '  +0.68%  '
$ws.Range("D49").Value = '0.06947'
$ws.Range("E49").Value = '  +1.32%  '
$ws.Range("D50").Value = '1.140'
$ws.Range("E50").Value = '  -0.47%  '
$ws.Range("D51").Value = '72.47'
$ws.Range("E51").Value = '  +0.77%  '
